# Apply the Feb 4 2024 cryptos-list refresh (prices + 1h volume deltas).
# Row 49/50 also swap identity (HuobiToken <-> MultiversX) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.650.73'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '2.279.72'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.51'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.501'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.56'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.04%  '
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').Value = '2.635.00'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').Value = '2.274.23'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.780'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').Value = '42.601.30'
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.06'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.61%  '
$ws.Range('D20').Value = '0.0₃0897'
$ws.Range('E20').Value = '  -1.69%  '
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.17'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.33%  '
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.85%  '
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.97'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.66'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.63%  '
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0691'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('E40').Value = '  -2.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.110'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.70'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.98%  '
$ws.Range('D43').Value = '2.000.71'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('E44').Value = '  -2.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.99'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.09'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.36%  '
$ws.Range('E48').Value = '  -1.87%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.71'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.19%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.85'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.39%  '
$ws.Range('D51').Value = '2.502.80'
$ws.Range('E51').Value = '  -0.96%  '
